$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing data rows (2-7) so the shared-string table is rebuilt
# cleanly in the exact order we (re)write the cells below.
$ws.Range("A2:L7").ClearContents()

# Column A (ΚΑΤΗΓΟΡΙΑ / price-list type)
$ws.Range("A2").Value = "Πελάτες Τιμή Πώλησης"
$ws.Range("A3").Value = "Πωλήσεις Έκπτωση 1"
$ws.Range("A4").Value = "Πωλήσεις Έκπτωση 1"
$ws.Range("A5").Value = "Πελάτες Τιμή Πώλησης"
$ws.Range("A6").Value = "Πελάτες Τιμή Πώλησης"
$ws.Range("A7").Value = "Πελάτες Τιμή Πώλησης"

# Column E (ΠΕΡΙΓΡΑΦΗ / description)
$ws.Range("E2").Value = "Βερύκοκα® Ελληνικά (Ζυγιζόμενο) /Kgr"
$ws.Range("E3").Value = "Nivea® Sun Spray Protect & Bronze {20} 200ml"
$ws.Range("E4").Value = "Παπαγάλος® Ελληνικός Καφές Κουπάτος 143gr"
$ws.Range("E5").Value = "Βερύκοκα® Ελληνικά (Ζυγιζόμενο) /Kgr"
$ws.Range("E6").Value = "Nutella® Πραλίνα Βάζο 400gr"
$ws.Range("E7").Value = "Νεκταρίνια® ΝΑΟΥΣΑΣ  (Ζυγιζόμενο) /Kgr"

# Column F (ΚΩΔΙΚΟΣ / code)
$ws.Range("F2").Value = "0253"
$ws.Range("F3").Value = "4005808859634"
$ws.Range("F4").Value = "5201219486417"
$ws.Range("F5").Value = "0253"
$ws.Range("F6").Value = "80135876"
$ws.Range("F7").Value = "0214"

# Column J (BRAND)
$ws.Range("J2").Value = "Βερύκοκα"
$ws.Range("J3").Value = "Nivea"
$ws.Range("J4").Value = "Παπαγάλος"
$ws.Range("J5").Value = "Βερύκοκα"
$ws.Range("J6").Value = "Nutella"
$ws.Range("J7").Value = "Νεκταρίνια"

# Dates (ΕΝΑΡΞΗ / ΛΗΞΗ) - unchanged, re-applied defensively
$ws.Range("C2").Value = 44028
$ws.Range("D2").Value = 44043
$ws.Range("C3").Value = 44028
$ws.Range("D3").Value = 44043
$ws.Range("C4").Value = 44028
$ws.Range("D4").Value = 44043
$ws.Range("C5").Value = 44028
$ws.Range("D5").Value = 44043
$ws.Range("C6").Value = 44028
$ws.Range("D6").Value = 44043
$ws.Range("C7").Value = 44028
$ws.Range("D7").Value = 44043

# Numeric columns G, H, I, K, L
$ws.Range("G2").Value = 1.95
$ws.Range("H2").Value = 1.65
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 3.84
$ws.Range("G3").Value = 8.9
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 50
$ws.Range("K3").Value = 4
$ws.Range("L3").Value = 14.34
$ws.Range("G4").Value = 3.7
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 50
$ws.Range("K4").Value = 13
$ws.Range("L4").Value = 21.23
$ws.Range("G5").Value = 1.95
$ws.Range("H5").Value = 1.65
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 37.15
$ws.Range("L5").Value = 54.27
$ws.Range("G6").Value = 3.78
$ws.Range("H6").Value = 2.95
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 42
$ws.Range("L6").Value = 110.36
$ws.Range("G7").Value = 0.85
$ws.Range("H7").Value = 0.85
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 177.943
$ws.Range("L7").Value = 134.03

# Totals row
$ws.Range("K8").Value = 277.093
$ws.Range("L8").Value = 338.0700000000001
